$wb = $excel.ActiveWorkbook

# ALC!row43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 2000
$ws.Cells.Item(43, 9).Value = 2000
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 2000
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = -1931
$ws.Cells.Item(43, 14).ClearContents()

# ALC!row75
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(75, 8).Value = 52650
$ws.Cells.Item(75, 10).Value = 52650
$ws.Cells.Item(75, 12).Value = 52650
$ws.Cells.Item(75, 14).Value = -54522

# ALC!row78
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(78, 8).Value = 52650
$ws.Cells.Item(78, 10).Value = 52650
$ws.Cells.Item(78, 12).Value = 157950
$ws.Cells.Item(78, 14).Value = -167310

# ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 2916.5
$ws.Cells.Item(100, 9).Value = 2199.8
$ws.Cells.Item(100, 11).Value = 2199.8
$ws.Cells.Item(100, 13).Value = -1658.8

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 2094.6667
$ws.Cells.Item(125, 10).Value = 2094.6667
$ws.Cells.Item(125, 12).Value = 18852.0003
$ws.Cells.Item(125, 14).Value = -23772.0003

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1471.1428
$ws.Cells.Item(45, 10).Value = 1316.6666
$ws.Cells.Item(45, 12).Value = 1316.6666
$ws.Cells.Item(45, 14).Value = -2070.6666

# ARM!row76
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 21578.8
$ws.Cells.Item(76, 10).Value = 21578.8
$ws.Cells.Item(76, 12).Value = 21578.8
$ws.Cells.Item(76, 14).Value = -22254.8

# ARM!row79
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 21578.8
$ws.Cells.Item(79, 10).Value = 21578.8
$ws.Cells.Item(79, 12).Value = 21578.8
$ws.Cells.Item(79, 14).Value = -23918.8

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 142858780
$ws.Cells.Item(102, 9).Value = 142858780
$ws.Cells.Item(102, 11).Value = 142858780
$ws.Cells.Item(102, 13).Value = -142857158

# BSM!row20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 14).ClearContents()

# BSM!row99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1013.5714
$ws.Cells.Item(99, 9).Value = 941.2
$ws.Cells.Item(99, 11).Value = 941.2
$ws.Cells.Item(99, 13).Value = 556.8

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5766.6665
$ws.Cells.Item(31, 9).Value = 5803.6665
$ws.Cells.Item(31, 11).Value = 5803.6665
$ws.Cells.Item(31, 13).Value = -5508.6665

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 5766.6665
$ws.Cells.Item(34, 9).Value = 5803.6665
$ws.Cells.Item(34, 11).Value = 5803.6665
$ws.Cells.Item(34, 13).Value = -5601.6665

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 328
$ws.Cells.Item(107, 9).Value = 235.5
$ws.Cells.Item(107, 11).Value = 235.5
$ws.Cells.Item(107, 13).Value = 1684.5

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 5665.1816
$ws.Cells.Item(132, 10).Value = 6500
$ws.Cells.Item(132, 12).Value = 19500
$ws.Cells.Item(132, 14).Value = -24560

# CUL!row33
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 96.125
$ws.Cells.Item(33, 10).Value = 125
$ws.Cells.Item(33, 12).Value = 750
$ws.Cells.Item(33, 14).Value = -1316

# CUL!row70
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 3739.25
$ws.Cells.Item(70, 9).Value = 2986
$ws.Cells.Item(70, 11).Value = 8958
$ws.Cells.Item(70, 13).Value = -8643

# CUL!row73
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(73, 8).Value = 3739.25
$ws.Cells.Item(73, 9).Value = 2986
$ws.Cells.Item(73, 11).Value = 8958
$ws.Cells.Item(73, 13).Value = -7866

# CUL!row75
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 220
$ws.Cells.Item(75, 10).Value = 225
$ws.Cells.Item(75, 12).Value = 675
$ws.Cells.Item(75, 14).Value = -2671

# CUL!row76
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(76, 8).Value = 6999.3335
$ws.Cells.Item(76, 9).Value = 5000
$ws.Cells.Item(76, 11).Value = 15000
$ws.Cells.Item(76, 13).Value = -14617

# CUL!row78
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(78, 8).Value = 220
$ws.Cells.Item(78, 10).Value = 225
$ws.Cells.Item(78, 12).Value = 2025
$ws.Cells.Item(78, 14).Value = -12009

# CUL!row79
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(79, 8).Value = 6999.3335
$ws.Cells.Item(79, 9).Value = 5000
$ws.Cells.Item(79, 11).Value = 15000
$ws.Cells.Item(79, 13).Value = -13674

# CUL!row80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 2000
$ws.Cells.Item(80, 9).Value = 1000
$ws.Cells.Item(80, 10).Value = 3000
$ws.Cells.Item(80, 11).Value = 3000
$ws.Cells.Item(80, 12).Value = 9000
$ws.Cells.Item(80, 13).Value = -2064
$ws.Cells.Item(80, 14).Value = -10872

# CUL!row83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value = 2000
$ws.Cells.Item(83, 9).Value = 1000
$ws.Cells.Item(83, 10).Value = 3000
$ws.Cells.Item(83, 11).Value = 9000
$ws.Cells.Item(83, 12).Value = 27000
$ws.Cells.Item(83, 13).Value = -4320
$ws.Cells.Item(83, 14).Value = -36360

# CUL!row98
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 608.5
$ws.Cells.Item(98, 10).Value = 625.7143
$ws.Cells.Item(98, 12).Value = 1877.1429
$ws.Cells.Item(98, 14).Value = -4873.1429

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 44995
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 44995
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 44995
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(70, 14).Value = -45535

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 44995
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 44995
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 44995
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(73, 14).Value = -46867

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 522.4
$ws.Cells.Item(102, 9).Value = 522.4
$ws.Cells.Item(102, 11).Value = 522.4
$ws.Cells.Item(102, 13).Value = 1099.6

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 10000
$ws.Cells.Item(16, 10).Value = 10000
$ws.Cells.Item(16, 12).Value = 10000
$ws.Cells.Item(16, 14).Value = -10340

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4040
$ws.Cells.Item(40, 9).Value = 1995.3334
$ws.Cells.Item(40, 11).Value = 1995.3334
$ws.Cells.Item(40, 13).Value = -1859.3334

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 575
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 13).ClearContents()

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1475
$ws.Cells.Item(61, 9).Value = 1475
$ws.Cells.Item(61, 11).Value = 1475
$ws.Cells.Item(61, 13).Value = -1273

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1475
$ws.Cells.Item(113, 9).Value = 1475
$ws.Cells.Item(113, 11).Value = 1475
$ws.Cells.Item(113, 13).Value = 695

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 37045940
$ws.Cells.Item(122, 9).Value = 66678924
$ws.Cells.Item(122, 10).Value = 4711.25
$ws.Cells.Item(122, 11).Value = 200036772
$ws.Cells.Item(122, 12).Value = 14133.75
$ws.Cells.Item(122, 13).Value = -200034322
$ws.Cells.Item(122, 14).Value = -19033.75

# WVR!row5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 7501940
$ws.Cells.Item(5, 10).Value = 7501940
$ws.Cells.Item(5, 12).Value = 7501940
$ws.Cells.Item(5, 14).Value = -7502164

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1897.4
$ws.Cells.Item(107, 10).Value = 2631.3333
$ws.Cells.Item(107, 12).Value = 7893.999899999999
$ws.Cells.Item(107, 14).Value = -11733.9999
